$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 9: fill in the two cells that were missing (X9, Y9) ---
$ws.Range("X9").Value = -1.2099989999999963
$ws.Range("Y9").Value = "Down"

# --- new row 10 (repeat-trader run appended) ---
$ws.Range("A10").Value = 42653.87903935185
$ws.Range("B10").Value = -9
$ws.Range("C10").Value = "Sell"
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = "Random"
$ws.Range("Q10").Value = 35.550971360736582
$ws.Range("R10").Value = -24.44
$ws.Range("S10").Value = -0.1153
$ws.Range("S10").NumberFormat = $ws.Range("S9").NumberFormat
$ws.Range("T10").Value = -0.047
$ws.Range("T10").NumberFormat = $ws.Range("T9").NumberFormat
$ws.Range("U10").Value = 6.45
$ws.Range("V10").Value = 1.88
$ws.Range("W10").Value = 1

# --- columns re-bestfit themselves a touch wider once the new row is in ---
$ws.Columns.Item(1).ColumnWidth = 14.5
$ws.Columns.Item(2).ColumnWidth = 7.6666666666666667
$ws.Columns.Item(3).ColumnWidth = 8.3333333333333333
$ws.Columns.Item(4).ColumnWidth = 11.3333333333333333
$ws.Columns.Item(5).ColumnWidth = 8.6666666666666667
$ws.Columns.Item(6).ColumnWidth = 11.3333333333333333
$ws.Columns.Item(7).ColumnWidth = 18.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 19.5
$ws.Columns.Item(10).ColumnWidth = 19.8333333333333333
$ws.Columns.Item(11).ColumnWidth = 9.5
$ws.Columns.Item(12).ColumnWidth = 13.5
$ws.Columns.Item(13).ColumnWidth = 13.8333333333333333
